$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet1 (25% Appreciation)
$ws1.Range("B2").Value = 3.6183943499999902
$ws1.Range("E2").Value = 361839.43499999901
$ws1.Range("B3").Value = 0.48308507912495602
$ws1.Range("E3").Value = 223107.74
$ws1.Range("B4").Value = 1.38285813792866
$ws1.Range("E4").Value = 947184.77500000002
$ws1.Range("B5").Value = 0.63053481623893903
$ws1.Range("E5").Value = 1047716.79625
$ws1.Range("B6").Value = 0.25738127612057898
$ws1.Range("E6").Value = 718616.73499999999
$ws1.Range("B7").Value = 0.35475937604982
$ws1.Range("E7").Value = 1245435.5215
$ws1.Range("B8").Value = 0.323234962186014
$ws1.Range("E8").Value = 1537332.68449999
$ws1.Range("B9").Value = 0.39419411588539999
$ws1.Range("E9").Value = 2480827.8899999899
$ws1.Range("B10").Value = 0.085777645630829194
$ws1.Range("E10").Value = 752634.05574999901
$ws1.Range("B11").Value = 0.058391231411115997
$ws1.Range("E11").Value = 579269.96600000001
$ws1.Range("B12").Value = 0.41095130466015201
$ws1.Range("E12").Value = 4315406.1294999998
$ws1.Range("B13").Value = 0.116563714058132
$ws1.Range("E13").Value = 1727057.1652499901
$ws1.Range("B14").Value = 0.131036198158442
$ws1.Range("E14").Value = 2167794.6004999899
$ws1.Range("B15").Value = 0.17976861291837501
$ws1.Range("E15").Value = 3410888.9730000002
$ws1.Range("B16").Value = 0.15224263048349701
$ws1.Range("E16").Value = 3407899.9219999998
$ws1.Range("B17").Value = 0.078045735004012601
$ws1.Range("E17").Value = 2015960.7749999899
$ws1.Range("B18").Value = 0.0795720738362415
$ws1.Range("E18").Value = 2215801.0580000002
$ws1.Range("B19").Value = 0.084681919201807093
$ws1.Range("E19").Value = 2545730.48
$ws1.Range("B20").Value = 0.033957224579987101
$ws1.Range("E20").Value = 1107277.09949999

# Sheet2 (50% Appreciation)
$ws2.Range("B3").Value = 0.66144180100000005
$ws2.Range("E3").Value = 192920.41499999899
$ws2.Range("B4").Value = 0.92118750299999996
$ws2.Range("E4").Value = 446395.41
$ws2.Range("B5").Value = 0.61716333899999998
$ws2.Range("E5").Value = 589750.37874999898
$ws2.Range("B6").Value = 0.27883087499999998
$ws2.Range("E6").Value = 437493
$ws2.Range("B7").Value = 0.41300601300000001
$ws2.Range("E7").Value = 828704.57099999895
$ws2.Range("B8").Value = 0.37460496599999998
$ws2.Range("E8").Value = 1062088.9724999999
$ws2.Range("B9").Value = 0.42866700099999999
$ws2.Range("E9").Value = 1670649.4449999901
$ws2.Range("B10").Value = 0.107586485
$ws2.Range("E10").Value = 599037.50124999904
$ws2.Range("B11").Value = 0.081455154000000002
$ws2.Range("E11").Value = 528330.25049999903
$ws2.Range("B12").Value = 0.37907299300000002
$ws2.Range("E12").Value = 2659473.4075000002
$ws2.Range("B13").Value = 0.096299396999999995
$ws2.Range("E13").Value = 931716.16024999996
$ws2.Range("B14").Value = 0.14480696700000001
$ws2.Range("E14").Value = 1535955.70999999
$ws2.Range("B15").Value = 0.20766116200000001
$ws2.Range("E15").Value = 2557944.091
$ws2.Range("B16").Value = 0.168964279
$ws2.Range("E16").Value = 2513481.9324999899
$ws2.Range("B17").Value = 0.075947081
$ws2.Range("E17").Value = 1323166.344
$ws2.Range("B18").Value = 0.10225379
$ws2.Range("E18").Value = 1916786.237
$ws2.Range("B19").Value = 0.095660239999999994
$ws2.Range("E19").Value = 1976547.8599999901
$ws2.Range("B20").Value = 0.041981499999999998
$ws2.Range("E20").Value = 950407.21500000299

# Sheet3 (75% Appreciation)
$ws3.Range("B3").Value = 0.59493041144237402
$ws3.Range("E3").Value = 120667.19500000001
$ws3.Range("B4").Value = 0.56146645247341898
$ws3.Range("E4").Value = 181630.424999999
$ws3.Range("B5").Value = 0.42544374335795798
$ws3.Range("E5").Value = 223325.35499999899
$ws3.Range("B6").Value = 0.37743068997359602
$ws3.Range("E6").Value = 282412.02500000002
$ws3.Range("B7").Value = 0.48863710951711498
$ws3.Range("E7").Value = 503619.08
$ws3.Range("B8").Value = 0.36419634564571302
$ws3.Range("E8").Value = 558779.1
$ws3.Range("B9").Value = 0.38487222034280899
$ws3.Range("E9").Value = 805560.23
$ws3.Range("B10").Value = 0.17883225156996599
$ws3.Range("E10").Value = 518366.58724999899
$ws3.Range("B11").Value = 0.1218766922108
$ws3.Range("E11").Value = 446271.80949999898
$ws3.Range("B12").Value = 0.29000947116618597
$ws3.Range("E12").Value = 1191341.04999999
$ws3.Range("B13").Value = 0.072225492514558898
$ws3.Range("E13").Value = 382743.06999999902
$ws3.Range("B14").Value = 0.13874159154784799
$ws3.Range("E14").Value = 788332.85099999898
$ws3.Range("B15").Value = 0.19499020082483101
$ws3.Range("E15").Value = 1278717.5404999999
$ws3.Range("B16").Value = 0.179476945666351
$ws3.Range("E16").Value = 1406484.17499999
$ws3.Range("B17").Value = 0.13506998322011601
$ws3.Range("E17").Value = 1252226.0847499899
$ws3.Range("B18").Value = 0.152712883947942
$ws3.Range("E18").Value = 1607023.4645
$ws3.Range("B19").Value = 0.142367894437082
$ws3.Range("E19").Value = 1726949.94
$ws3.Range("B20").Value = 0.047307926743142499
$ws3.Range("E20").Value = 655552.67099999997

# Sheet4 (100% Appreciation)
$ws4.Range("B3").Value = 0.496576809614165
$ws4.Range("E3").Value = 78881.22
$ws4.Range("B4").Value = 0.27033682567737699
$ws4.Range("E4").Value = 64267.5
$ws4.Range("B5").Value = 0.31979248685106998
$ws4.Range("E5").Value = 101373.80499999999
$ws4.Range("B6").Value = 0.45814422380342201
$ws4.Range("E6").Value = 191674.95
$ws4.Range("B7").Value = 0.41043141786693299
$ws4.Range("E7").Value = 250382.64499999999
$ws4.Range("B8").Value = 0.36899248089130798
$ws4.Range("E8").Value = 317492.24
$ws4.Range("B9").Value = 0.35892346884946003
$ws4.Range("E9").Value = 422783.97499999899
$ws4.Range("B10").Value = 0.26235243020351001
$ws4.Range("E10").Value = 419949.19374999899
$ws4.Range("B11").Value = 0.167926577316208
$ws4.Range("E11").Value = 362457.31999999902
$ws4.Range("B12").Value = 0.27364329227853801
$ws4.Range("E12").Value = 689823.18500000006
$ws4.Range("B13").Value = 0.088275281142178599
$ws4.Range("E13").Value = 283426.13999999902
$ws4.Range("B14").Value = 0.075595795596146903
$ws4.Range("E14").Value = 264141.840499999
$ws4.Range("B15").Value = 0.18310055651704599
$ws4.Range("E15").Value = 688142.4
$ws4.Range("B16").Value = 0.209911511840492
$ws4.Range("E16").Value = 933354.37749999994
$ws4.Range("B17").Value = 0.194194725827222
$ws4.Range("E17").Value = 1044723.47074999
$ws4.Range("B18").Value = 0.19492697189901001
$ws4.Range("E18").Value = 1252307.57075
$ws4.Range("B19").Value = 0.11055851190174799
$ws4.Range("E19").Value = 848735.99
$ws4.Range("B20").Value = 0.12164736588753999
$ws4.Range("E20").Value = 1037109.44249999

# Set column D width on sheet1 ("25% Appreciation") to a best-fit-style width
$ws1.Columns.Item(4).ColumnWidth = 9.83

# Replicate the final navigation/selection state recorded in the sheetViews:
# user moved through each sheet and ended with the "75% Appreciation" tab active.
$ws1.Activate()
$ws1.Range("H21").Select()

$ws2.Activate()
$ws2.Range("F30").Select()

$ws4.Activate()
$ws4.Range("E36").Select()

$ws3.Activate()
$ws3.Range("F36").Select()
